# edit.ps1 - reproduce the "Reverted naming in architecture figure" commit
#
# Main content edit: on the (only) slide, three small rectangle labels in the
# flex-routing architecture diagram get a "Main" prefix:
#   Control  -> MainControl
#   Deparser -> MainDeparser
#   Parser   -> MainParser
#
# Incidental edit also present in the target OOXML: the cached text of the
# auto-updating "datetimeFigureOut" date fields on the slide master / slide
# layouts (12/11/2022 -> 12/13/2022) and the notes master (11/12/2022 ->
# 13/12/2022) was refreshed (e.g. by a re-save around the 12/12/2022 meeting
# mentioned in the commit message). We refresh every such placeholder we can
# reach through the object model.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Rename the three architecture-diagram labels on slide 1.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -eq "Control") {
                $shp.TextFrame.TextRange.Text = "MainControl"
            } elseif ($t -eq "Deparser") {
                $shp.TextFrame.TextRange.Text = "MainDeparser"
            } elseif ($t -eq "Parser") {
                $shp.TextFrame.TextRange.Text = "MainParser"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Refresh the cached "today" text of the date placeholders.
# ---------------------------------------------------------------------------

# 2a. Slide master(s) + every layout belonging to them.
for ($d = 1; $d -le $p.Designs.Count; $d++) {
    $master = $p.Designs.Item($d).SlideMaster

    for ($j = 1; $j -le $master.Shapes.Count; $j++) {
        $shp = $master.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq "12/11/2022") {
                    $shp.TextFrame.TextRange.Text = "12/13/2022"
                }
            }
        }
    }

    for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
        $lay = $master.CustomLayouts.Item($i)
        for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
            $shp = $lay.Shapes.Item($j)
            if ($shp.HasTextFrame) {
                if ($shp.TextFrame.HasText) {
                    if ($shp.TextFrame.TextRange.Text -eq "12/11/2022") {
                        $shp.TextFrame.TextRange.Text = "12/13/2022"
                    }
                }
            }
        }
    }
}

# 2b. Notes master.
#
# NOTE: deliberately *not* touched here. This runtime's NotesMaster shape
# access aliases onto the slide master's shape collection by numeric shape
# Id (the notes master's "Date Placeholder 2" happens to share Id=3 with the
# slide master's "Text Placeholder 2"), so writing through
# `$p.NotesMaster.Shapes.Item(...)` silently corrupts unrelated slide-master
# text instead of updating the notes master. Left unchanged to avoid that
# corruption; every other occurrence of the date-field text is still
# refreshed above/below.
